$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new values look numeric (e.g. "598.38") need to be forced to
# text explicitly, otherwise Excel will silently convert them to real numbers
# (losing formatting like trailing zeros) -- the source data models these as text.
$numericLookingCells = @("D5","D6","D8","D10","D13","D18","D20","D21","D23","D25","D26","D27","D28","D29","D31","D33","D34","D35","D36","D37","D40","D42","D43","D44","D45","D46","D47","D48","D49","D50")
foreach ($addr in $numericLookingCells) {
    $ws.Range($addr).NumberFormat = "@"
}

# Apply updated values
$ws.Range('D2').Value = '65.833.90'
$ws.Range('D3').Value = '2.657.54'
$ws.Range('E3').Value = '  -0.16%  '
$ws.Range('E4').Value = '  -0.02%  '
$ws.Range('D5').Value = '598.38'
$ws.Range('E5').Value = '  -0.39%  '
$ws.Range('D6').Value = '157.40'
$ws.Range('E6').Value = '  +0.33%  '
$ws.Range('E7').Value = '  -0.03%  '
$ws.Range('D8').Value = '0.630'
$ws.Range('E8').Value = '  +1.73%  '
$ws.Range('E9').Value = '  +2.13%  '
$ws.Range('D10').Value = '0.397'
$ws.Range('E10').Value = '  -1.15%  '
$ws.Range('E11').Value = '  -0.87%  '
$ws.Range('E12').Value = '  +1.24%  '
$ws.Range('D13').Value = '28.62'
$ws.Range('E13').Value = '  -2.37%  '
$ws.Range('E14').Value = '  +0.35%  '
$ws.Range('D15').Value = '3.132.95'
$ws.Range('E15').Value = '  -0.13%  '
$ws.Range('D16').Value = '65.622.81'
$ws.Range('E16').Value = '  +0.29%  '
$ws.Range('D17').Value = '2.642.32'
$ws.Range('E17').Value = '  -1.14%  '
$ws.Range('D18').Value = '12.58'
$ws.Range('E18').Value = '  -0.54%  '
$ws.Range('E19').Value = '  -1.36%  '
$ws.Range('D20').Value = '7.46'
$ws.Range('E20').Value = '  -1.68%  '
$ws.Range('D21').Value = '349.88'
$ws.Range('E21').Value = '  -0.35%  '
$ws.Range('E22').Value = '  +0.13%  '
$ws.Range('D23').Value = '69.21'
$ws.Range('E23').Value = '  -0.19%  '
$ws.Range('E24').Value = '  +2.62%  '
$ws.Range('D25').Value = '1.74'
$ws.Range('E25').Value = '  +7.12%  '
$ws.Range('D26').Value = '9.64'
$ws.Range('E26').Value = '  -0.71%  '
$ws.Range('D27').Value = '1.60'
$ws.Range('E27').Value = '  +0.88%  '
$ws.Range('B28').Value = 'Bittensor'
$ws.Range('C28').Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Range('D28').Value = '553.45'
$ws.Range('E28').Value = '  +3.96%  '
$ws.Range('B29').Value = 'Kaspa'
$ws.Range('C29').Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range('D29').Value = '0.164'
$ws.Range('E29').Value = '  -2.09%  '
$ws.Range('E30').Value = '  +0.00%  '
$ws.Range('D31').Value = '7.92'
$ws.Range('E31').Value = '  -1.67%  '
$ws.Range('E32').Value = '  +0.58%  '
$ws.Range('D33').Value = '1.77'
$ws.Range('E33').Value = '  +1.16%  '
$ws.Range('D34').Value = '6.49'
$ws.Range('E34').Value = '  -0.33%  '
$ws.Range('D35').Value = '5.43'
$ws.Range('E35').Value = '  -0.74%  '
$ws.Range('D36').Value = '0.419'
$ws.Range('E36').Value = '  -0.93%  '
$ws.Range('D37').Value = '20.35'
$ws.Range('E37').Value = '  -0.02%  '
$ws.Range('E38').Value = '  -0.06%  '
$ws.Range('E39').Value = '  +0.22%  '
$ws.Range('D40').Value = '154.78'
$ws.Range('E40').Value = '  -2.74%  '
$ws.Range('E41').Value = '  +0.01%  '
$ws.Range('D42').Value = '161.73'
$ws.Range('E42').Value = '  -1.59%  '
$ws.Range('D43').Value = '4.06'
$ws.Range('E43').Value = '  -0.23%  '
$ws.Range('D44').Value = '2.29'
$ws.Range('E44').Value = '  +0.60%  '
$ws.Range('D45').Value = '0.0605'
$ws.Range('D46').Value = '22.62'
$ws.Range('E46').Value = '  -1.24%  '
$ws.Range('D47').Value = '0.638'
$ws.Range('E47').Value = '  -0.56%  '
$ws.Range('D48').Value = '0.0256'
$ws.Range('E48').Value = '  -0.89%  '
$ws.Range('D49').Value = '0.100'
$ws.Range('E49').Value = '  -1.24%  '
$ws.Range('D50').Value = '19.73'
$ws.Range('E50').Value = '  -2.15%  '
$ws.Range('E51').Value = '  +7.15%  '
